$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.825.71"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").Value = "'1.813.45"
$ws.Range("E3").Value = "  +1.00%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'309.40"
$ws.Range("E5").Value = "  +0.41%  "

$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").Value = "'0.4331"
$ws.Range("E7").Value = "  +2.83%  "

$ws.Range("E8").Value = "  +3.16%  "

$ws.Range("D9").Value = "'0.07252"
$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").Value = "'0.8667"
$ws.Range("E10").Value = "  +3.04%  "

$ws.Range("D11").Value = "'20.95"
$ws.Range("E11").Value = "  +3.57%  "

$ws.Range("D12").Value = "'1.939.91"
$ws.Range("E12").Value = "  +7.11%  "

$ws.Range("D13").Value = "'6.675"
$ws.Range("E13").Value = "  +4.97%  "

$ws.Range("D14").Value = "'5.365"
$ws.Range("E14").Value = "  +1.70%  "

$ws.Range("D15").Value = "'0.06918"
$ws.Range("E15").Value = "  +2.08%  "

$ws.Range("D16").Value = "'1.009"
$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("D17").Value = "'80.55"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").Value = "'0.000008934"
$ws.Range("E18").Value = "  +2.34%  "

$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").Value = "'15.23"
$ws.Range("E20").Value = "  +1.51%  "

$ws.Range("D21").Value = "'26.866.21"
$ws.Range("E21").Value = "  -0.83%  "

$ws.Range("D22").Value = "'5.217"
$ws.Range("E22").Value = "  +2.81%  "

$ws.Range("D23").Value = "'11.19"
$ws.Range("E23").Value = "  +1.32%  "

$ws.Range("D24").Value = "'2.174.38"
$ws.Range("E24").Value = "  +6.07%  "

$ws.Range("D25").Value = "'153.82"
$ws.Range("E25").Value = "  +0.57%  "

$ws.Range("D26").Value = "'1.872"
$ws.Range("E26").Value = "  -4.10%  "

$ws.Range("D27").Value = "'18.30"
$ws.Range("E27").Value = "  +0.88%  "

$ws.Range("D28").Value = "'5.221"
$ws.Range("E28").Value = "  +4.30%  "

$ws.Range("D29").Value = "'1.912"
$ws.Range("E29").Value = "  +16.00%  "

$ws.Range("D30").Value = "'115.36"
$ws.Range("E30").Value = "  +1.65%  "

$ws.Range("D31").Value = "'0.08931"
$ws.Range("E31").Value = "  -0.58%  "

$ws.Range("D32").Value = "'0.7569"
$ws.Range("E32").Value = "  +4.33%  "

$ws.Range("D33").Value = "'1.172"
$ws.Range("E33").Value = "  +7.62%  "

$ws.Range("D34").Value = "'4.436"
$ws.Range("E34").Value = "  +2.36%  "

$ws.Range("D35").Value = "'2.805"
$ws.Range("E35").Value = "  -1.97%  "

$ws.Range("D36").Value = "'1.007"
$ws.Range("E36").Value = "  +0.30%  "

$ws.Range("D37").Value = "'1.132"
$ws.Range("E37").Value = "  +5.09%  "

$ws.Range("D38").Value = "'0.05230"
$ws.Range("E38").Value = "  +1.74%  "

$ws.Range("D39").Value = "'0.01928"
$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("D40").Value = "'0.5084"
$ws.Range("E40").Value = "  +2.38%  "

$ws.Range("E41").Value = "  +1.36%  "

$ws.Range("D42").Value = "'2.678"
$ws.Range("E42").Value = "  +2.36%  "

$ws.Range("D43").Value = "'6.556"
$ws.Range("E43").Value = "  +10.41%  "

$ws.Range("D44").Value = "'8.286"
$ws.Range("E44").Value = "  +2.89%  "

$ws.Range("D45").Value = "'106.74"
$ws.Range("E45").Value = "  +1.73%  "

$ws.Range("D46").Value = "'10.40"
$ws.Range("E46").Value = "  +1.47%  "

$ws.Range("D47").Value = "'1.004"
$ws.Range("E47").Value = "  +0.08%  "

$ws.Range("E48").Value = "  +3.66%  "

$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "'0.4560"
$ws.Range("E49").Value = "  +0.79%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06281"
$ws.Range("E50").Value = "  -0.44%  "

$ws.Range("D51").Value = "'1.809"
$ws.Range("E51").Value = "  +5.42%  "

# Reset number format/style on touched D/E cells so no stray quotePrefix style lingers on them
$ws.Range("D2:E51").Style = "Normal"

